# ZBP_07_testovani.xlsx update:
#   Sheet "data"   -> add new column AH ("24. 8. 2021") with per-row % values,
#                     update the footer note date (A59).
#   Sheet "pocetR" -> add new column AG ("24. 8. 2021") with per-row sample
#                     sizes, add a blank trailing cell (AG21), update the
#                     footer note date (A21).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "data": new column AH
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Bring over the header formatting (style s="1") from the previous last
# column (AG) before writing the new header label.
$wsData.Range("AG1").Copy()
$wsData.Range("AH1").PasteSpecial(-4122)
$wsData.Range("AH1").Value = "24. 8. 2021"

$dataValues = "0.06,0.03,0.22,0.03,0.06,0.32,0.05,0.02,0.28,0.07000000000000001,0.02,0.1,0.07000000000000001,0.03,0.24,0.04,0.02,0.2,0.06,0.03,0.27,0.05,0.08,0.27,0.06,0.015,0.09,0.03,0.04,0.31,0.07000000000000001,0.03,0.27,0.05,0.06,0.35,0.01,0.11,0.14,0.05,0.03,0.26,0.03,0.05,0.24,0.05,0.05,0.26,0.09,0.02,0.32,0.09,0.015,0.21,0.05,0.08,0.19".Split(",")

$row = 2
foreach ($v in $dataValues) {
    $wsData.Cells.Item($row, 34).Value = [double]$v
    $row = $row + 1
}

# Footer note (row 59, column A) - bump the "aktualizace" date.
$wsData.Range("A59").Value = "Život během pandemie, Testování, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# ---------------------------------------------------------------------------
# Sheet "pocetR": new column AG
# ---------------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

$wsPocet.Range("AF1").Copy()
$wsPocet.Range("AG1").PasteSpecial(-4122)
$wsPocet.Range("AG1").Value = "24. 8. 2021"

$pocetValues = "1901,452,709,740,937,964,917,156,586,242,671,76,65,261,348,102,308,175,96".Split(",")

$row = 2
foreach ($v in $pocetValues) {
    $wsPocet.Cells.Item($row, 33).Value = [double]$v
    $row = $row + 1
}

# Footer note (row 21, column A) - bump the "aktualizace" date.
$wsPocet.Range("A21").Value = "Život během pandemie, Testování, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"

# Trailing blank cell matching the rest of the footer row (B21:AF21 are all
# empty-string cells already; extend that through the new AG column).
$wsPocet.Range("AG21").Value = ""
